# Update attendance report: Y2_B2526_Respiratory_session_analysis.xlsx
#
# The underlying roster sizes for several groups changed (new students
# added/removed), a couple of sessions got newly recorded, and the list of
# "Recorded By" staff for a few already-recorded sessions was refreshed.
# This script reproduces those data edits on the "Session Analysis Results"
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) "Students" column (H): denominator (roster size) refreshed per group
# ---------------------------------------------------------------------

# Group A1 (rows 2-16): 202 -> 216
foreach ($r in 2..16) {
    $ws.Range("H$r").Value = "0/216"
}

# Group A2 (rows 17-31): 215 -> 217
foreach ($r in 17..31) {
    $ws.Range("H$r").Value = "0/217"
}

# Group B3 (rows 92-106): 225 -> 224
foreach ($r in 92..106) {
    $ws.Range("H$r").Value = "0/224"
}

# Group B4 (rows 107-121): 225 -> 226
foreach ($r in 107..121) {
    $ws.Range("H$r").Value = "0/226"
}

# ---------------------------------------------------------------------
# 2) Individually recorded sessions: attendance numerator + "Recorded By"
# ---------------------------------------------------------------------

$ws.Range("G32").Value = "mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("H32").Value = "157/220"

$ws.Range("G47").Value = "mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("H47").Value = "113/225"

$ws.Range("G62").Value = "mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("H62").Value = "94/154"

$ws.Range("G77").Value = "mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("H77").Value = "153/224"

$ws.Range("H82").Value = "123/224"

# Row 67 (HISTOLOGY B1 #1) flips from "Not Recorded" to "Recorded": copy the
# green "Recorded" row formatting (matches the style used by e.g. row 32)
# across the whole row, then fill in the recorder + attendance + status.
$ws.Range("I32").Copy()
$ws.Range("A67:I67").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("G67").Value = "Sara_nabil@med.asu.edu.eg"
$ws.Range("H67").Value = "23/154"
$ws.Range("I67").Value = "Recorded"

# ---------------------------------------------------------------------
# 3) "Class Statistics" summary box (K3:L10)
# ---------------------------------------------------------------------

$ws.Range("L4").Value = 1706   # Total Students
$ws.Range("L6").Value = 6      # Recorded Sessions
$ws.Range("L7").Value = 2      # Missing Sessions

# Percentage values are stored as literal text in this sheet (not numeric
# percentages). Force text entry with a leading apostrophe, then restore the
# original (non-percent) number formatting by copying it from an untouched
# neighboring cell that uses the same base style.
$ws.Range("L9").Value = "'5.0%"    # Coverage %
$ws.Range("L10").Value = "'53.5%"  # Average Attendance %

$ws.Range("L8").Copy()
$ws.Range("L9:L10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4) "Group Statistics" table (K14:S22)
# ---------------------------------------------------------------------

$ws.Range("M15").Value = 216   # A1 Students
$ws.Range("M16").Value = 217   # A2 Students

$ws.Range("S17").Value = "'71.4%"  # A2 Avg Attendance %
$ws.Range("S18").Value = "'50.2%"  # A3 Avg Attendance %

$ws.Range("O19").Value = 2     # B1 Recorded
$ws.Range("P19").Value = 0     # B1 Missing
$ws.Range("R19").Value = "'13.3%"  # B1 Coverage %
$ws.Range("S19").Value = "'38.0%"  # B1 Avg Attendance %

$ws.Range("S20").Value = "'61.6%"  # B2 Avg Attendance %

$ws.Range("M21").Value = 224   # B3 Students
$ws.Range("M22").Value = 226   # B4 Students

$ws.Range("L8").Copy()
$ws.Range("S17:S18").PasteSpecial(-4122)
$ws.Range("R19:S19").PasteSpecial(-4122)
$ws.Range("S20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "Applied attendance report updates."
